$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 762
$wsExhibit.Range("F4").Value = 1519
$wsExhibit.Range("F5").Value = 232
$wsExhibit.Range("F6").Value = 97
$wsExhibit.Range("F7").Value = 158
$wsExhibit.Range("F8").Value = 6281
$wsExhibit.Range("F9").Value = 71
$wsExhibit.Range("F10").Value = 407
$wsExhibit.Range("F11").Value = 118
$wsExhibit.Range("F12").Value = 5274
$wsExhibit.Range("F13").Value = 30
$wsExhibit.Range("F14").Value = 180
$wsExhibit.Range("F15").Value = 1189
$wsExhibit.Range("F16").Value = 1189
$wsExhibit.Range("F17").Value = 60
$wsExhibit.Range("F19").Value = 72
$wsExhibit.Range("F20").Value = 12
$wsExhibit.Range("F22").Value = 29
$wsExhibit.Range("F23").Value = 3770
$wsExhibit.Range("F24").Value = 163

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 87

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 87
$wsAll.Range("F4").Value = 762
$wsAll.Range("F5").Value = 1519
$wsAll.Range("F6").Value = 232
$wsAll.Range("F7").Value = 97
$wsAll.Range("F8").Value = 158
$wsAll.Range("F9").Value = 6281
$wsAll.Range("F10").Value = 71
$wsAll.Range("F11").Value = 407
$wsAll.Range("F12").Value = 118
$wsAll.Range("F13").Value = 5274
$wsAll.Range("F14").Value = 30
$wsAll.Range("F15").Value = 180
$wsAll.Range("F16").Value = 1189
$wsAll.Range("F17").Value = 1189
$wsAll.Range("F18").Value = 60
$wsAll.Range("F20").Value = 72
$wsAll.Range("F21").Value = 12
$wsAll.Range("F23").Value = 29
$wsAll.Range("F24").Value = 3770
$wsAll.Range("F26").Value = 163
